# Auto-generated edit script: update market health data values
$wb = $excel.ActiveWorkbook

# --- Sheet: Metadata ---
$ws = $wb.Worksheets.Item("Metadata")
$ws.Cells.Item(2, 1).Value = '30 Oct 2025, 11:34 AM'

# --- Sheet: Top Losers ---
$ws = $wb.Worksheets.Item("Top Losers")
$ws.Cells.Item(2, 3).Value = -18.5527
$ws.Cells.Item(2, 4).Value = -17.2921
$ws.Cells.Item(2, 5).Value = -0.264
$ws.Cells.Item(3, 3).Value = -8.735
$ws.Cells.Item(3, 4).Value = -5.267
$ws.Cells.Item(3, 5).Value = 7.0789
$ws.Cells.Item(4, 3).Value = -5.5556
$ws.Cells.Item(4, 4).Value = -8.1081
$ws.Cells.Item(4, 5).Value = 8.7331
$ws.Cells.Item(12, 2).Value = 'RAJRATAN'
$ws.Cells.Item(12, 3).Value = -4.1125
$ws.Cells.Item(12, 4).Value = -3.327
$ws.Cells.Item(12, 5).Value = 21.5882
$ws.Cells.Item(13, 2).Value = 'KALAMANDIR'
$ws.Cells.Item(13, 3).Value = -3.9601
$ws.Cells.Item(13, 4).Value = -2.2841
$ws.Cells.Item(13, 5).Value = 21.0099
$ws.Cells.Item(14, 3).Value = -3.6088
$ws.Cells.Item(14, 4).Value = 17.3683
$ws.Cells.Item(14, 5).Value = 18.1356
$ws.Cells.Item(15, 3).Value = -3.5723
$ws.Cells.Item(15, 4).Value = -1.2511
$ws.Cells.Item(15, 5).Value = 1.2832
$ws.Cells.Item(16, 2).Value = 'ATHERENERG'
$ws.Cells.Item(16, 3).Value = -3.5472
$ws.Cells.Item(16, 4).Value = -3.5608
$ws.Cells.Item(16, 5).Value = 20.4509
$ws.Cells.Item(17, 2).Value = 'TVSHLTD'
$ws.Cells.Item(17, 3).Value = -3.4813
$ws.Cells.Item(17, 4).Value = -2.2385
$ws.Cells.Item(17, 5).Value = 16.0266
$ws.Cells.Item(18, 2).Value = 'YATRA'
$ws.Cells.Item(18, 3).Value = -3.3091
$ws.Cells.Item(18, 4).Value = -6.0604
$ws.Cells.Item(18, 5).Value = 3.818
$ws.Cells.Item(19, 2).Value = 'UTIAMC'
$ws.Cells.Item(19, 3).Value = -3.2895
$ws.Cells.Item(19, 4).Value = -7.8714
$ws.Cells.Item(19, 5).Value = -5.314
$ws.Cells.Item(20, 2).Value = 'SARDAEN'
$ws.Cells.Item(20, 3).Value = -3.1851
$ws.Cells.Item(20, 4).Value = 0.3431
$ws.Cells.Item(20, 5).Value = 0.3059
$ws.Cells.Item(21, 2).Value = 'GOKULAGRO'
$ws.Cells.Item(21, 3).Value = -3.1742
$ws.Cells.Item(21, 4).Value = 4.3773
$ws.Cells.Item(21, 5).Value = -13.9867
$ws.Cells.Item(22, 3).Value = -3.1343
$ws.Cells.Item(22, 4).Value = 3.7313
$ws.Cells.Item(22, 5).Value = 11.7183
$ws.Cells.Item(23, 3).Value = -3.1056
$ws.Cells.Item(23, 4).Value = 6.8493
$ws.Cells.Item(23, 5).Value = 22.0896
$ws.Cells.Item(25, 2).Value = 'KHAICHEM'
$ws.Cells.Item(25, 3).Value = -3.0363
$ws.Cells.Item(25, 4).Value = -9.5528
$ws.Cells.Item(25, 5).Value = -7.5407
$ws.Cells.Item(26, 2).Value = 'UBL'
$ws.Cells.Item(26, 3).Value = -2.9071
$ws.Cells.Item(26, 4).Value = -2.3126
$ws.Cells.Item(26, 5).Value = -0.7774
$ws.Cells.Item(27, 2).Value = 'JSL'
$ws.Cells.Item(27, 3).Value = -2.8307
$ws.Cells.Item(27, 4).Value = -2.6308
$ws.Cells.Item(27, 5).Value = 5.8556
$ws.Cells.Item(28, 2).Value = 'SAIL'
$ws.Cells.Item(28, 3).Value = -2.7748
$ws.Cells.Item(28, 4).Value = 5.5538
$ws.Cells.Item(28, 5).Value = 1.6136
$ws.Cells.Item(29, 2).Value = 'PRECWIRE'
$ws.Cells.Item(29, 3).Value = -2.7362
$ws.Cells.Item(29, 4).Value = 9.2823
$ws.Cells.Item(29, 5).Value = 19.9128
$ws.Cells.Item(30, 2).Value = 'QUESS'
$ws.Cells.Item(30, 3).Value = -2.7059
$ws.Cells.Item(30, 4).Value = 4.4425
$ws.Cells.Item(30, 5).Value = -3.967
$ws.Cells.Item(31, 2).Value = 'SANDHAR'
$ws.Cells.Item(31, 3).Value = -2.6608
$ws.Cells.Item(31, 4).Value = 1.0551
$ws.Cells.Item(31, 5).Value = 18.4219
$ws.Cells.Item(32, 2).Value = 'EPACKPEB'
$ws.Cells.Item(32, 3).Value = -2.6419
$ws.Cells.Item(32, 4).Value = -2.5766
$ws.Cells.Item(32, 5).Value = 'N/A'
$ws.Cells.Item(34, 2).Value = 'INDUSTOWER'
$ws.Cells.Item(34, 3).Value = -2.5587
$ws.Cells.Item(34, 4).Value = 2.6967
$ws.Cells.Item(34, 5).Value = 8.2823
$ws.Cells.Item(35, 3).Value = -2.5543
$ws.Cells.Item(35, 4).Value = -2.9935
$ws.Cells.Item(35, 5).Value = -4.1673
$ws.Cells.Item(36, 2).Value = 'MEGASOFT'
$ws.Cells.Item(36, 3).Value = -2.5388
$ws.Cells.Item(36, 4).Value = 12.8199
$ws.Cells.Item(36, 5).Value = 30.1372
$ws.Cells.Item(37, 3).Value = -2.5341
$ws.Cells.Item(37, 4).Value = 5.2644
$ws.Cells.Item(37, 5).Value = 6.8718
$ws.Cells.Item(38, 2).Value = 'SURAJEST'
$ws.Cells.Item(38, 3).Value = -2.5071
$ws.Cells.Item(38, 4).Value = 6.4805
$ws.Cells.Item(38, 5).Value = 4.4779
$ws.Cells.Item(40, 2).Value = 'DCMSRIND'
$ws.Cells.Item(40, 3).Value = -2.4115
$ws.Cells.Item(40, 4).Value = -1.1492
$ws.Cells.Item(40, 5).Value = 4.7101
$ws.Cells.Item(41, 2).Value = 'AEROFLEX'
$ws.Cells.Item(41, 3).Value = -2.2833
$ws.Cells.Item(41, 4).Value = 5.0485
$ws.Cells.Item(41, 5).Value = 3.8969
$ws.Cells.Item(42, 2).Value = 'SINDHUTRAD'
$ws.Cells.Item(42, 3).Value = -2.2829
$ws.Cells.Item(42, 4).Value = -1.1954
$ws.Cells.Item(42, 5).Value = -15.1504
$ws.Cells.Item(43, 2).Value = 'CPPLUS'
$ws.Cells.Item(43, 3).Value = -2.2536
$ws.Cells.Item(43, 4).Value = -2.9079
$ws.Cells.Item(43, 5).Value = 2.2436
$ws.Cells.Item(44, 2).Value = 'GABRIEL'
$ws.Cells.Item(44, 3).Value = -2.2415
$ws.Cells.Item(44, 4).Value = 1.8917
$ws.Cells.Item(44, 5).Value = 6.687
$ws.Cells.Item(45, 2).Value = 'HFCL'
$ws.Cells.Item(45, 3).Value = -2.2294
$ws.Cells.Item(45, 4).Value = -2.9714
$ws.Cells.Item(45, 5).Value = 3.4989
$ws.Cells.Item(46, 2).Value = 'VGUARD'
$ws.Cells.Item(46, 3).Value = -2.1958
$ws.Cells.Item(46, 4).Value = 0.5239
$ws.Cells.Item(46, 5).Value = -0.3595
$ws.Cells.Item(47, 2).Value = 'NEWGEN'
$ws.Cells.Item(47, 3).Value = -2.1843
$ws.Cells.Item(47, 4).Value = 9.1011
$ws.Cells.Item(47, 5).Value = 9.4751
$ws.Cells.Item(48, 2).Value = 'HONASA'
$ws.Cells.Item(48, 3).Value = -2.1766
$ws.Cells.Item(48, 4).Value = -0.7781
$ws.Cells.Item(48, 5).Value = -2.403
$ws.Cells.Item(49, 2).Value = 'SOUTHBANK'
$ws.Cells.Item(49, 3).Value = -2.1478
$ws.Cells.Item(49, 4).Value = 0.4198
$ws.Cells.Item(49, 5).Value = 32.1934
$ws.Cells.Item(50, 2).Value = 'HMT'
$ws.Cells.Item(50, 3).Value = -2.1376
$ws.Cells.Item(50, 4).Value = -2.6284
$ws.Cells.Item(50, 5).Value = -6.0247
$ws.Cells.Item(51, 2).Value = 'POCL'
$ws.Cells.Item(51, 3).Value = -2.1375
$ws.Cells.Item(51, 4).Value = 3.0446
$ws.Cells.Item(51, 5).Value = 23.641
$ws.Cells.Item(52, 2).Value = 'JTEKTINDIA'
$ws.Cells.Item(52, 3).Value = -2.1265
$ws.Cells.Item(52, 4).Value = 4.1621
$ws.Cells.Item(52, 5).Value = -1.8386
$ws.Cells.Item(53, 2).Value = 'TVSELECT'
$ws.Cells.Item(53, 3).Value = -2.1142
$ws.Cells.Item(53, 4).Value = -3.0674
$ws.Cells.Item(53, 5).Value = -5.0477
$ws.Cells.Item(54, 2).Value = 'VINCOFE'
$ws.Cells.Item(54, 3).Value = -2.0919
$ws.Cells.Item(54, 4).Value = 12.5533
$ws.Cells.Item(54, 5).Value = 10.8997
$ws.Cells.Item(55, 2).Value = 'CHAMBLFERT'
$ws.Cells.Item(55, 3).Value = -2.069
$ws.Cells.Item(55, 4).Value = -0.584
$ws.Cells.Item(55, 5).Value = -5.0122
$ws.Cells.Item(56, 2).Value = 'HCG'
$ws.Cells.Item(56, 3).Value = -2.0482
$ws.Cells.Item(56, 4).Value = 0.1713
$ws.Cells.Item(56, 5).Value = 18.203
$ws.Cells.Item(58, 2).Value = 'GRWRHITECH'
$ws.Cells.Item(58, 3).Value = -2.0267
$ws.Cells.Item(58, 4).Value = -5.6327
$ws.Cells.Item(58, 5).Value = 19.2109
$ws.Cells.Item(59, 2).Value = 'MOLDTKPAC'
$ws.Cells.Item(59, 3).Value = -1.9986
$ws.Cells.Item(59, 4).Value = -3.2199
$ws.Cells.Item(59, 5).Value = -1.5419
$ws.Cells.Item(60, 2).Value = 'IDBI'
$ws.Cells.Item(60, 3).Value = -1.9821
$ws.Cells.Item(60, 4).Value = 6.3451
$ws.Cells.Item(60, 5).Value = 9.2649
$ws.Cells.Item(61, 2).Value = 'CROMPTON'
$ws.Cells.Item(61, 3).Value = -1.9584
$ws.Cells.Item(61, 4).Value = -2.7271
$ws.Cells.Item(61, 5).Value = -2.0426
$ws.Cells.Item(62, 2).Value = 'AMBER'
$ws.Cells.Item(62, 3).Value = -1.9549
$ws.Cells.Item(62, 4).Value = -2.0609
$ws.Cells.Item(62, 5).Value = 0.7541
$ws.Cells.Item(63, 2).Value = 'ARIHANTCAP'
$ws.Cells.Item(63, 3).Value = -1.9328
$ws.Cells.Item(63, 4).Value = 5.0558
$ws.Cells.Item(63, 5).Value = -3.8012
$ws.Cells.Item(64, 2).Value = 'SPLPETRO'
$ws.Cells.Item(64, 3).Value = -1.9153
$ws.Cells.Item(64, 4).Value = -6.8432
$ws.Cells.Item(64, 5).Value = -9.5355
$ws.Cells.Item(65, 2).Value = 'APARINDS'
$ws.Cells.Item(65, 3).Value = -1.9129
$ws.Cells.Item(65, 4).Value = 6.269
$ws.Cells.Item(65, 5).Value = 13.3765
$ws.Cells.Item(66, 2).Value = 'SUMMITSEC'
$ws.Cells.Item(66, 3).Value = -1.9041
$ws.Cells.Item(66, 4).Value = -3.5203
$ws.Cells.Item(66, 5).Value = 3.9912
$ws.Cells.Item(67, 2).Value = 'APOLLOPIPE'
$ws.Cells.Item(67, 3).Value = -1.9038
$ws.Cells.Item(67, 4).Value = -3.6145
$ws.Cells.Item(67, 5).Value = -8.7087
$ws.Cells.Item(68, 2).Value = 'RAJRILTD'
$ws.Cells.Item(68, 3).Value = -1.9
$ws.Cells.Item(68, 4).Value = -4.962
$ws.Cells.Item(68, 5).Value = -14.5683
$ws.Cells.Item(69, 2).Value = 'MSPL'
$ws.Cells.Item(69, 3).Value = -1.895
$ws.Cells.Item(69, 4).Value = -0.7375
$ws.Cells.Item(69, 5).Value = -7.8082
$ws.Cells.Item(70, 2).Value = 'HITECHGEAR'
$ws.Cells.Item(70, 3).Value = -1.8947
$ws.Cells.Item(70, 4).Value = -0.7678
$ws.Cells.Item(70, 5).Value = 7.8426
$ws.Cells.Item(71, 2).Value = 'KPITTECH'
$ws.Cells.Item(71, 3).Value = -1.8926
$ws.Cells.Item(71, 4).Value = -0.6108
$ws.Cells.Item(71, 5).Value = 6.7232
$ws.Cells.Item(72, 2).Value = 'DCBBANK'
$ws.Cells.Item(72, 3).Value = -1.8794
$ws.Cells.Item(72, 4).Value = -0.9957
$ws.Cells.Item(72, 5).Value = 23.1655
$ws.Cells.Item(73, 2).Value = 'BHARATWIRE'
$ws.Cells.Item(73, 3).Value = -1.8712
$ws.Cells.Item(73, 4).Value = 20.5351
$ws.Cells.Item(73, 5).Value = 21.5795
$ws.Cells.Item(74, 2).Value = 'TIIL'
$ws.Cells.Item(74, 3).Value = -1.8533
$ws.Cells.Item(74, 4).Value = -1.5577
$ws.Cells.Item(74, 5).Value = 11.1766
$ws.Cells.Item(75, 2).Value = 'RAYMONDREL'
$ws.Cells.Item(75, 3).Value = -1.8462
$ws.Cells.Item(75, 4).Value = -3.4431
$ws.Cells.Item(75, 5).Value = 11.2467
$ws.Cells.Item(76, 2).Value = 'PROSTARM'
$ws.Cells.Item(76, 3).Value = -1.8443
$ws.Cells.Item(76, 4).Value = -1.1846
$ws.Cells.Item(76, 5).Value = -9.9216

# --- Sheet: 1 Month Performance ---
$ws = $wb.Worksheets.Item("1 Month Performance")
$ws.Cells.Item(3, 3).Value = 81.8182
$ws.Cells.Item(4, 3).Value = 78.86
$ws.Cells.Item(5, 3).Value = 67.0098
$ws.Cells.Item(6, 3).Value = 61.353
$ws.Cells.Item(7, 3).Value = 58.8745
$ws.Cells.Item(9, 3).Value = 53.6263
$ws.Cells.Item(10, 3).Value = 47.6867
$ws.Cells.Item(13, 3).Value = 39.985
$ws.Cells.Item(14, 3).Value = 39.2391
$ws.Cells.Item(16, 2).Value = 'BHARATSE'
$ws.Cells.Item(16, 3).Value = 37.2885
$ws.Cells.Item(17, 2).Value = 'SHAREINDIA'
$ws.Cells.Item(17, 3).Value = 37.2266
$ws.Cells.Item(18, 3).Value = 35.9859
$ws.Cells.Item(19, 3).Value = 35.3607
$ws.Cells.Item(21, 2).Value = 'RAMAPHO'
$ws.Cells.Item(21, 3).Value = 34.1011
$ws.Cells.Item(22, 2).Value = 'MEGASOFT'
$ws.Cells.Item(22, 3).Value = 34.0059
$ws.Cells.Item(24, 3).Value = 32.5139
$ws.Cells.Item(25, 3).Value = 31.3443
$ws.Cells.Item(26, 3).Value = 29.6147
$ws.Cells.Item(27, 3).Value = 29.6083
$ws.Cells.Item(30, 3).Value = 27.1415
$ws.Cells.Item(31, 2).Value = 'SAGILITY'
$ws.Cells.Item(31, 3).Value = 26.7382
$ws.Cells.Item(33, 2).Value = 'EMKAY'
$ws.Cells.Item(33, 3).Value = 26.6455
$ws.Cells.Item(34, 3).Value = 26.2074
$ws.Cells.Item(35, 3).Value = 25.9652
$ws.Cells.Item(36, 3).Value = 25.8046
$ws.Cells.Item(37, 3).Value = 25.4928
$ws.Cells.Item(38, 3).Value = 24.7282
$ws.Cells.Item(39, 3).Value = 24.3943
$ws.Cells.Item(41, 3).Value = 24.0198
$ws.Cells.Item(42, 3).Value = 23.8611
$ws.Cells.Item(43, 2).Value = 'CARTRADE'
$ws.Cells.Item(43, 3).Value = 23.8362
$ws.Cells.Item(44, 2).Value = 'TDPOWERSYS'
$ws.Cells.Item(44, 3).Value = 23.6878
$ws.Cells.Item(45, 2).Value = 'KERNEX'
$ws.Cells.Item(45, 3).Value = 23.5882
$ws.Cells.Item(46, 3).Value = 23.2712
$ws.Cells.Item(47, 2).Value = 'TATVA'
$ws.Cells.Item(47, 3).Value = 23.2354
$ws.Cells.Item(48, 2).Value = 'LORDSCHLO'
$ws.Cells.Item(48, 3).Value = 22.9389
$ws.Cells.Item(51, 3).Value = 22.4429
$ws.Cells.Item(52, 3).Value = 22.2256
$ws.Cells.Item(53, 2).Value = 'SURYODAY'
$ws.Cells.Item(53, 3).Value = 21.9087
$ws.Cells.Item(54, 3).Value = 21.9051
$ws.Cells.Item(55, 2).Value = 'SCI'
$ws.Cells.Item(55, 3).Value = 21.8541
$ws.Cells.Item(56, 3).Value = 21.7359
$ws.Cells.Item(57, 2).Value = 'GUJTHEM'
$ws.Cells.Item(57, 3).Value = 21.3377
$ws.Cells.Item(58, 2).Value = 'PRIVISCL'
$ws.Cells.Item(58, 3).Value = 21.3028
$ws.Cells.Item(61, 3).Value = 20.7169
$ws.Cells.Item(63, 3).Value = 20.4174
$ws.Cells.Item(64, 3).Value = 20.0545
$ws.Cells.Item(65, 3).Value = 19.8
$ws.Cells.Item(66, 3).Value = 19.7684
$ws.Cells.Item(67, 2).Value = 'MCX'
$ws.Cells.Item(67, 3).Value = 19.6194
$ws.Cells.Item(68, 2).Value = 'HINDCOPPER'
$ws.Cells.Item(68, 3).Value = 19.5889
$ws.Cells.Item(69, 2).Value = 'WHEELS'
$ws.Cells.Item(69, 3).Value = 19.189
$ws.Cells.Item(70, 2).Value = 'BHAGERIA'
$ws.Cells.Item(70, 3).Value = 19.0221
$ws.Cells.Item(71, 2).Value = 'FIVESTAR'
$ws.Cells.Item(71, 3).Value = 18.8631
$ws.Cells.Item(72, 2).Value = 'HINDPETRO'
$ws.Cells.Item(72, 3).Value = 18.8557
$ws.Cells.Item(73, 2).Value = 'PRECWIRE'
$ws.Cells.Item(73, 3).Value = 18.748
$ws.Cells.Item(74, 2).Value = 'ACUTAAS'
$ws.Cells.Item(74, 3).Value = 18.6626
$ws.Cells.Item(75, 2).Value = 'REPRO'
$ws.Cells.Item(75, 3).Value = 18.6236
$ws.Cells.Item(76, 2).Value = 'ETHOSLTD'
$ws.Cells.Item(76, 3).Value = 18.42

# --- Sheet: distance from Dma50 ---
$ws = $wb.Worksheets.Item("distance from Dma50")
$ws.Cells.Item(2, 3).Value = 10.0487
$ws.Cells.Item(3, 3).Value = 7.5855
$ws.Cells.Item(4, 3).Value = 6.7276
$ws.Cells.Item(5, 3).Value = 5.5217
$ws.Cells.Item(6, 3).Value = 5.3548
$ws.Cells.Item(7, 3).Value = 5.2868
$ws.Cells.Item(8, 3).Value = 4.5768
$ws.Cells.Item(9, 3).Value = 4.5049
$ws.Cells.Item(10, 3).Value = 4.0087
$ws.Cells.Item(11, 3).Value = 3.8062
$ws.Cells.Item(12, 3).Value = 3.5724
$ws.Cells.Item(13, 3).Value = 3.5262
$ws.Cells.Item(14, 3).Value = 3.2642
$ws.Cells.Item(15, 3).Value = 3.2172
$ws.Cells.Item(16, 3).Value = 3.1428
$ws.Cells.Item(17, 3).Value = 2.9835
$ws.Cells.Item(18, 3).Value = 2.889
$ws.Cells.Item(19, 3).Value = 2.8205
$ws.Cells.Item(20, 3).Value = 2.5072
$ws.Cells.Item(21, 3).Value = 2.417
$ws.Cells.Item(22, 3).Value = 1.5363
$ws.Cells.Item(23, 3).Value = 1.4691
$ws.Cells.Item(24, 3).Value = 1.4429
$ws.Cells.Item(25, 3).Value = 1.1302
$ws.Cells.Item(26, 3).Value = 1.1058
$ws.Cells.Item(27, 3).Value = 0.9101
$ws.Cells.Item(28, 3).Value = 0.8253
$ws.Cells.Item(29, 3).Value = 0.4297
$ws.Cells.Item(30, 3).Value = -2.2098

